$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'black pants tight',
    'black pants women',
    'black pastel',
    'black pastels',
    'black plus leggings',
    'black plus size leggings 3x',
    'black plus size stockings',
    'black pocket hose',
    'black purple',
    'black rainbow',
    'black rainbows',
    'black red compression',
    'black riding pants for women',
    'black riding tights',
    'black running leggings for women',
    'black running pants',
    'black running pants women',
    'black running tights',
    'black seam thigh highs',
    'black seamless leggings for women',
    'black sheer knee highs',
    'black sheer leggings',
    'black sheer panty hose',
    'black sheer tight',
    'black shimmer tights women',
    'black short leggings',
    'black short pants',
    'black short tights',
    'black shorts for women tight',
    'black skirt leggings',
    'black soccer pants women',
    'black soccer shorts',
    'black spandex',
    'black sport',
    'black sport leggings for women',
    'black sport pants women',
    'black sport tape',
    'black sports leggings',
    'black sports leggings women',
    'black squad game',
    'black stage tape',
    'black stocking',
    'black stocking plus size',
    'black stockings',
    'black stockings knee high',
    'black stockings plus size',
    'black stripe shorts women',
    'black struggle red scare',
    'black summer leggings for women',
    'black support',
    'black support knee highs for women',
    'black tear away pants',
    'black thigh',
    'black thigh high',
    'black thigh highs',
    'black thigh highs for women',
    'black tight',
    'black tight pants women',
    'black tight shorts',
    'black tight small',
    'black tight stockings',
    'black tights',
    'black tights control top women',
    'black tights for women',
    'black tights for women with designs',
    'black tights leggings',
    'black tights opaque',
    'black tights pack',
    'black tights with seam',
    'black tights women',
    'black tights women opaque',
    'black tights womens',
    'black tighys',
    'black top for leggings',
    'black tops to wear with leggings',
    'black trousers for women',
    'black tummy control tights',
    'black underarmour pants',
    'black warmup suit',
    'black web tights',
    'black winter over the knee skirts for women',
    'black winter stockings for women',
    'black woman pants',
    'black woman tights',
    'black women clothes',
    'black women compression shorts',
    'black women ring',
    'black women tights',
    'black womens capri leggings',
    'black womens pants',
    'black womens tights',
    'black workout capris for women',
    'black workout clothes',
    'black workout leggings for women',
    'black workout leggings women',
    'black x',
    'black x tights',
    'black xs for women',
    'black yoga',
    'black yoga capri'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
